$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the raw "Inputs" data (columns H:O) for rows 54-92. These values
#    were previously missing (rows 54-68 had template formulas evaluating to
#    0, rows 69-92 did not exist at all). The B:F formulas already present
#    for rows 54-68 are shared formulas that will recompute automatically;
#    for the brand-new rows 69-92 we add the same formulas below.
# ---------------------------------------------------------------------------
$csv = @"
28434707,1712,32400,10242,28434082,1087,31775,9617
28434082,1492,40032,11984,28433499,909,39449,11401
28433499,3156,114258,49023,28444169,13826,124928,59693
28444169,809,26275,34570,28453045,9685,35151,43446
28453045,1801,50764,9881,28452919,1675,50638,9755
28452919,1989,75061,57292,28473717,22787,95859,78090
28473717,8637,136136,21016,28471156,6076,133575,18455
28471156,3220,60403,1924618,28472597,4661,61844,1926059
28472597,4045,150389,93813,28492680,24128,170472,113896
28492680,3758,150036,72380,28507778,18856,165134,87478
28507778,189,20644,8314,28507589,0,20455,8125
28507589,5118,76267,43207,28516920,14449,85598,52538
28516920,2168,39737,1947834,28526856,12104,49673,1957770
28526856,0,0,248,28526856,0,0,248
28526856,2373,47241,3999,28538368,13885,58753,15511
28538368,0,0,51655,28538368,0,0,51655
28538368,6329,102835,50318,28543409,11370,107876,55359
28543409,7771,164996,78851,28559270,23632,180857,94712
28559270,8113,143665,59019,28565619,14462,150014,65368
28565619,1749,17179,39647,28564659,789,16219,38687
28564659,2205,75770,1953589,28586663,24209,97774,1975593
28586663,4500,186729,78326,28584885,2722,184951,76548
28584885,2600,78071,57060,28597460,15175,90646,69635
28597460,4752,172101,91229,28620903,28195,195544,114672
28620903,3691,80989,1939895,28621693,4481,81779,1940685
28621693,2829,86139,2006861,28636112,17248,100558,2021280
28636112,3976,125061,23171,28643034,10898,131983,30093
28643034,83084,84240,1989214,28560051,101,1257,1906231
28560051,249655,379671,2140540,28322949,12553,142569,1903438
28322949,5300,83804,76275,28340157,22508,101012,93483
28340157,1944,95318,79763,28366982,28769,122143,106588
28366982,1608,70470,28509,28386696,21322,90184,48223
28386696,3536,70100,70135,28399679,16519,83083,83118
28399679,1360,5358,1910021,28399385,1066,5064,1909727
28399385,2118,134651,124484,28433781,36514,169047,158880
28433781,1827,89965,85277,28461215,29261,117399,112711
28461215,0,0,0,28461215,0,0,0
28461215,0,0,0,28461215,0,0,0
28461215,9367,392301,187205,28496308,44460,427394,222298
"@

$lines = $csv -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$rowCount = $lines.Count
$data = New-Object 'object[,]' $rowCount,8

for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $lines[$i].Trim() -split ","
    for ($j = 0; $j -lt 8; $j++) {
        $data[$i,$j] = [double]$parts[$j]
    }
}

$startRow = 54
$endRow = $startRow + $rowCount - 1
$ws.Range("H$startRow`:O$endRow").Value = $data

# ---------------------------------------------------------------------------
# 2. Column A dates for the new rows (69-92) -- rows 54-68 already had these.
# ---------------------------------------------------------------------------
for ($r = 69; $r -le 92; $r++) {
    $ws.Cells.Item($r, 1).Value = $r + 41602
    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item($r - 1, 1).Style
}

# ---------------------------------------------------------------------------
# 3. B:F formulas for the brand-new rows 69-92 (rows 54-68 already carry
#    these shared formulas and recalc automatically once H:O is populated).
# ---------------------------------------------------------------------------
$ws.Range("B69:B92").Formula = "=L69"
$ws.Range("C69:C92").Formula = "=M69"
$ws.Range("D69:D92").Formula = "=I69"
$ws.Range("E69:E92").Formula = "=N69-M69"
$ws.Range("F69:F92").Formula = "=O69-M69"

# ---------------------------------------------------------------------------
# 4. Summary rows 95-97: Average / Min / Max of C:F across the data rows.
# ---------------------------------------------------------------------------
$ws.Range("A95").Value = "Average"
$ws.Range("C95").Formula = "=AVERAGE(C3:C92)"
$ws.Range("D95:F95").Formula = "=AVERAGE(D3:D92)"
$ws.Range("C95:F95").NumberFormat = "0"

$ws.Range("A96").Value = "Min"
$ws.Range("C96").Formula = "=MIN(C3:C92)"
$ws.Range("D96:F96").Formula = "=MIN(D3:D92)"

$ws.Range("A97").Value = "Max"
$ws.Range("C97").Formula = "=MAX(C3:C92)"
$ws.Range("D97:F97").Formula = "=MAX(D3:D92)"

# ---------------------------------------------------------------------------
# 5. Misc formatting touch-up matching the author's edit: a couple of cells
#    in row 69 (I69:J69) picked up the black-font "Author Count" style.
# ---------------------------------------------------------------------------
$ws.Range("I69:J69").Font.Color = 0

# ---------------------------------------------------------------------------
# 6. View state: move the selection the same way the author left it.
# ---------------------------------------------------------------------------
$ws.Range("E104").Select() | Out-Null
